$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,("ECs", "Gas6", "Tyro3", "ECs", 3, 1, 19.23438166666667, 57.70314500000001, 0.1061888747949397, 0.1061888747949397, 2, 0.6666666666666666, 0.2333953333333333, 0.700186, 0.0679343338760815, 0.0679343338760815, 4.489214920552222, 40.40293428497, 0.007213870474244847, 0.007213870474244846)
  ,("ECs", "Gas6", "Tyro3", "FAPs", 3, 1, 19.23438166666667, 57.70314500000001, 0.1061888747949397, 0.1061888747949397, 3, 1, 2.855037666666667, 8.565113, 0.8310152534160438, 0.8310152534160438, 54.91488415337611, 494.233957380385, 0.08824457469768134, 0.08824457469768132)
  ,("ECs", "Gas6", "Tyro3", "sCs", 3, 1, 19.23438166666667, 57.70314500000001, 0.1061888747949397, 0.1061888747949397, 3, 1, 0.347169, 1.041507, 0.1010504127078748, 0.1010504127078748, 6.677581048835, 60.098229439515, 0.0107304296230135, 0.01073042962301349)
  ,("FAPs", "Gas6", "Tyro3", "ECs", 3, 1, 101.9328183333333, 305.798455, 0.5627491161960234, 0.5627491161960234, 2, 0.6666666666666666, 0.2333953333333333, 0.700186, 0.0679343338760815, 0.0679343338760815, 23.79064411251444, 214.11579701263, 0.03822998634813043, 0.03822998634813043)
  ,("FAPs", "Gas6", "Tyro3", "FAPs", 3, 1, 101.9328183333333, 305.798455, 0.5627491161960234, 0.5627491161960234, 3, 1, 2.855037666666667, 8.565113, 0.8310152534160438, 0.8310152534160438, 291.0220358111572, 2619.198322300415, 0.467653099405293, 0.467653099405293)
  ,("FAPs", "Gas6", "Tyro3", "sCs", 3, 1, 101.9328183333333, 305.798455, 0.5627491161960234, 0.5627491161960234, 3, 1, 0.347169, 1.041507, 0.1010504127078748, 0.1010504127078748, 35.387914607965, 318.491231471685, 0.05686603044259995, 0.05686603044259995)
  ,("M1", "Gas6", "Tyro3", "ECs", 3, 1, 33.990832, 101.972496, 0.1876560560134372, 0.1876560560134371, 2, 0.6666666666666666, 0.2333953333333333, 0.700186, 0.0679343338760815, 0.0679343338760815, 7.933301564917334, 71.399714084256, 0.01274828916308549, 0.01274828916308549)
  ,("M1", "Gas6", "Tyro3", "FAPs", 3, 1, 33.990832, 101.972496, 0.1876560560134372, 0.1876560560134371, 3, 1, 2.855037666666667, 8.565113, 0.8310152534160438, 0.8310152534160438, 97.04510568133868, 873.4059511320481, 0.1559450449430618, 0.1559450449430617)
  ,("M1", "Gas6", "Tyro3", "sCs", 3, 1, 33.990832, 101.972496, 0.1876560560134372, 0.1876560560134371, 3, 1, 0.347169, 1.041507, 0.1010504127078748, 0.1010504127078748, 11.800563154608, 106.205068391472, 0.01896272190728989, 0.01896272190728989)
  ,("M2", "Gas6", "Tyro3", "ECs", 3, 1, 24.872162, 74.61648600000001, 0.1373138446698593, 0.1373138446698593, 2, 0.6666666666666666, 0.2333953333333333, 0.700186, 0.0679343338760815, 0.0679343338760815, 5.805046540710667, 52.245418866396, 0.009328324569610618, 0.009328324569610615)
  ,("M2", "Gas6", "Tyro3", "FAPs", 3, 1, 24.872162, 74.61648600000001, 0.1373138446698593, 0.1373138446698593, 3, 1, 2.855037666666667, 8.565113, 0.8310152534160438, 0.8310152534160438, 71.01095936143534, 639.0986342529181, 0.1141098994258544, 0.1141098994258544)
  ,("M2", "Gas6", "Tyro3", "sCs", 3, 1, 24.872162, 74.61648600000001, 0.1373138446698593, 0.1373138446698593, 3, 1, 0.347169, 1.041507, 0.1010504127078748, 0.1010504127078748, 8.634843609378, 77.713592484402, 0.0138756206743943, 0.01387562067439429)
  ,("sCs", "Gas6", "Tyro3", "ECs", 3, 1, 1.103486, 3.310458, 0.006092108325740414, 0.006092108325740414, 2, 0.6666666666666666, 0.2333953333333333, 0.700186, 0.0679343338760815, 0.0679343338760815, 0.2575484827986667, 2.317936345188, 0.0004138633210101052, 0.0004138633210101051)
  ,("sCs", "Gas6", "Tyro3", "FAPs", 3, 1, 1.103486, 3.310458, 0.006092108325740414, 0.006092108325740414, 3, 1, 2.855037666666667, 8.565113, 0.8310152534160438, 0.8310152534160438, 3.150494094639333, 28.354446851754, 0.005062634944153161, 0.00506263494415316)
  ,("sCs", "Gas6", "Tyro3", "sCs", 3, 1, 1.103486, 3.310458, 0.006092108325740414, 0.006092108325740414, 3, 1, 0.347169, 1.041507, 0.1010504127078748, 0.1010504127078748, 0.383096131134, 3.447865180206, 0.0006156100605771489, 0.0006156100605771488)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = 2 + $i
  for ($j = 0; $j -lt 20; $j++) {
    $col = $j + 1
    $ws.Cells.Item($row, $col).Value = $data[$i][$j]
  }
}